# Update vm_pu.xlsx results for the "case with 380 kV" run.
# Slack/ext-grid voltage setpoint (column B) moves from 1.05 pu to 1.02 pu,
# and the resulting bus voltage magnitudes (columns C-F, I-M) are refreshed
# for every data row (rows 2-25). Columns G and N are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newResults = @{
    2 = @{ "B"=1.02; "C"=1.024525819036797; "D"=1.027687710516158; "E"=1.024933552443948; "F"=1.023009657280859; "I"=1.028540493957157; "J"=1.029699942307616; "K"=1.03050671189918; "L"=1.027760601293353; "M"=1.025842355591445 }
    3 = @{ "B"=1.02; "C"=1.026161437136498; "D"=1.028856640750343; "E"=1.026344343702474; "F"=1.025292674300836; "I"=1.028848906323447; "J"=1.030970959808858; "K"=1.031482444217269; "L"=1.028976949128002; "M"=1.027928137798961 }
    4 = @{ "B"=1.02; "C"=1.027216559348742; "D"=1.029610172946366; "E"=1.027254609729758; "F"=1.026766015625233; "I"=1.02904592541324; "J"=1.03178989909272; "K"=1.032110404993197; "L"=1.029760893682251; "M"=1.029273558623409 }
    5 = @{ "B"=1.02; "C"=1.027659374512954; "D"=1.029926286923337; "E"=1.027636672956733; "F"=1.027384496626049; "I"=1.029128146644078; "J"=1.03213335714089; "K"=1.032373594298924; "L"=1.03008972997923; "M"=1.029838188597965 }
    6 = @{ "B"=1.02; "C"=1.027733681114168; "D"=1.029979324790005; "E"=1.027700787484404; "F"=1.02748828964073; "I"=1.029141916539691; "J"=1.032190977296248; "K"=1.032417737981618; "L"=1.030144900329333; "M"=1.02993293542691 }
    7 = @{ "B"=1.02; "C"=1.027222479224694; "D"=1.02961439949967; "E"=1.027259717267628; "F"=1.026774283343722; "I"=1.02904702643206; "J"=1.031794491616827; "K"=1.032113924893461; "L"=1.029765290477369; "M"=1.029281107072734 }
    8 = @{ "B"=1.02; "C"=1.025079262460586; "D"=1.02808335076466; "E"=1.025410883861994; "F"=1.023782042605536; "I"=1.028645252316325; "J"=1.030130219620264; "K"=1.030837175534647; "L"=1.028172323415551; "M"=1.026548144470262 }
    9 = @{ "B"=1.02; "C"=1.021277184534424; "D"=1.025363220781516; "E"=1.022132456322145; "F"=1.018478039050548; "I"=1.0279176286784; "J"=1.027170236782255; "K"=1.028560892051068; "L"=1.02534094549041; "M"=1.021698855559106 }
    10 = @{ "B"=1.02; "C"=1.018724362017634; "D"=1.023534227281019; "E"=1.019932240805937; "F"=1.01491924672758; "I"=1.027419111115213; "J"=1.025177751421224; "K"=1.027024963942713; "L"=1.023436247308853; "M"=1.018441894768158 }
    11 = @{ "B"=1.02; "C"=1.017614445821065; "D"=1.022738417258138; "E"=1.018975885074752; "F"=1.013372452707878; "I"=1.027200009760929; "J"=1.024310259714223; "K"=1.026355390141623; "L"=1.022607268868115; "M"=1.017025518474922 }
    12 = @{ "B"=1.02; "C"=1.017201474107578; "D"=1.022442228136365; "E"=1.018620088630881; "F"=1.012796997185136; "I"=1.027118134505714; "J"=1.023987308178796; "K"=1.026105991807054; "L"=1.022298699342182; "M"=1.016498467029757 }
    13 = @{ "B"=1.02; "C"=1.017290089904134; "D"=1.022505788600579; "E"=1.018696433916071; "F"=1.012920475957933; "I"=1.027135719329091; "J"=1.024056615505821; "K"=1.026159519956526; "L"=1.022364918199222; "M"=1.016611564677435 }
    14 = @{ "B"=1.02; "C"=1.017580323830693; "D"=1.022713946277783; "E"=1.018946486410964; "F"=1.013324904082716; "I"=1.0271932519745; "J"=1.024283579350206; "K"=1.026334788917443; "L"=1.022581775758439; "M"=1.016981971722612 }
    15 = @{ "B"=1.02; "C"=1.017759053271298; "D"=1.022842120559513; "E"=1.019100476896283; "F"=1.013573964341206; "I"=1.027228634504375; "J"=1.024423322530948; "K"=1.026442686281114; "L"=1.022715302292909; "M"=1.017210065404292 }
    16 = @{ "B"=1.02; "C"=1.018797926479437; "D"=1.023586960536372; "E"=1.019995632829095; "F"=1.015021776983073; "I"=1.027433583506643; "J"=1.025235222964499; "K"=1.027069305401476; "L"=1.023491173509505; "M"=1.018535764023185 }
    17 = @{ "B"=1.02; "C"=1.019448358622434; "D"=1.024053141217993; "E"=1.020556153124788; "F"=1.015928371361886; "I"=1.027561272023931; "J"=1.025743228784279; "K"=1.027461152205052; "L"=1.023976714099404; "M"=1.019365687213086 }
    18 = @{ "B"=1.02; "C"=1.019827309191882; "D"=1.024324686499199; "E"=1.020882744204919; "F"=1.016456613926625; "I"=1.027635438273585; "J"=1.026039084434718; "K"=1.027689275776093; "L"=1.02425951440155; "M"=1.019849181822854 }
    19 = @{ "B"=1.02; "C"=1.019956448254305; "D"=1.024417214031584; "E"=1.020994044259679; "F"=1.016636637262416; "I"=1.027660674210939; "J"=1.026139886762494; "K"=1.027766986749718; "L"=1.024355873338586; "M"=1.020013942607701 }
    20 = @{ "B"=1.02; "C"=1.019378618593231; "D"=1.024003162805108; "E"=1.020496051009811; "F"=1.015831160387785; "I"=1.027547604589697; "J"=1.02568877179887; "K"=1.027419155730671; "L"=1.023924662422042; "M"=1.01927670509919 }
    21 = @{ "B"=1.02; "C"=1.017494876636564; "D"=1.022652665412852; "E"=1.018872867858495; "F"=1.013205835402476; "I"=1.027176323641874; "J"=1.02421676433522; "K"=1.026283195686443; "L"=1.022517934671761; "M"=1.016872922452826 }
    22 = @{ "B"=1.02; "C"=1.016306436949061; "D"=1.021800135179709; "E"=1.01784904112612; "F"=1.011549922551789; "I"=1.026940040039956; "J"=1.023287043731903; "K"=1.02556498133815; "L"=1.021629701051387; "M"=1.015356078569646 }
    23 = @{ "B"=1.02; "C"=1.016936842306824; "D"=1.02225240583533; "E"=1.018392105887928; "F"=1.012428263980413; "I"=1.027065569593559; "J"=1.023780310609135; "K"=1.025946102627334; "L"=1.022100932486619; "M"=1.016160717149333 }
    24 = @{ "B"=1.02; "C"=1.019410132455771; "D"=1.024025747041799; "E"=1.020523209652602; "F"=1.015875087567518; "I"=1.027553781279861; "J"=1.025713379973826; "K"=1.027438133470701; "L"=1.023948183592812; "M"=1.019316914094581 }
    25 = @{ "B"=1.02; "C"=1.022263230495427; "D"=1.026069137816638; "E"=1.022982523262792; "F"=1.019853138259466; "I"=1.028108088081285; "J"=1.027938782195218; "K"=1.029152565777022; "L"=1.026075886678193; "M"=1.02295664205883 }
}

foreach ($row in $newResults.Keys) {
    $rowData = $newResults[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
